$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 28
    3  = 27
    4  = 27
    5  = 27
    6  = 27
    7  = 28
    8  = 28
    9  = 27
    10 = 28
    11 = 28
    12 = 30
    13 = 29
    14 = 29
    15 = 24
    16 = 26
    17 = 24
    18 = 25
    19 = 23
    20 = 25
    21 = 23
    22 = 23
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
